# "Start of day 3 part 2"
# The existing day3 part 1 grid (B11:G16) is relocated ten columns to the
# right and five rows down, to L16:Q21, making room above/left of it for
# the part-2 work. Each row is moved individually (rather than as one big
# B11:G16 -> L16 block) so that the sparse rows (11 and 16, which only had
# a single populated cell) don't get padded with extra blank cells in the
# destination - that keeps the moved sheet shaped exactly like the
# original.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the grid, row by row (Cut clears the source and moves the values).
$ws.Range("G11").Cut($ws.Range("Q16"))
$ws.Range("B12:G12").Cut($ws.Range("L17"))
$ws.Range("B13:G13").Cut($ws.Range("L18"))
$ws.Range("B14:G14").Cut($ws.Range("L19"))
$ws.Range("B15:G15").Cut($ws.Range("L20"))
$ws.Range("B16:G16").Cut($ws.Range("L21"))

# Cut only moves values, so re-create the formulas at their new locations
# (references shifted by the same +10 columns / +5 rows offset).
$ws.Range("Q16").Formula = "=Q17+P17"
$ws.Range("Q17").Formula = "=Q18+P18+P17"
$ws.Range("Q18").Formula = "=Q19+P19+P18+P17"
$ws.Range("Q19").Formula = "=Q20+P20+P19+P18"
$ws.Range("Q20").Formula = "=Q21+P21+P20+P19"
$ws.Range("O21").Formula = "=N21+N20+O20"
$ws.Range("P21").Formula = "=O21+O20+P20"
$ws.Range("Q21").Formula = "=P21+P20"

# Recalculate so the cached formula results are correct.
$excel.Calculate()

# Match the author's new cursor position/selection.
[void]$ws.Range("P16").Select()
